$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 27 column B holds a date-like string ("2025-12-30"); set text format first
# so Excel does not auto-convert it to a date serial (matches original inlineStr cells).
$ws.Range("B27").NumberFormat = "@"

$ws.Cells.Item(3, 7).Value = 34
$ws.Cells.Item(3, 10).Value = 5
$ws.Cells.Item(4, 6).Value = 2.08
$ws.Cells.Item(4, 9).Value = 4.5
$ws.Cells.Item(4, 11).Value = 3.85
$ws.Cells.Item(4, 16).Value = 1.83
$ws.Cells.Item(4, 17).Value = 1.97
$ws.Cells.Item(5, 9).Value = 6.4
$ws.Cells.Item(5, 16).Value = 2.56
$ws.Cells.Item(5, 17).Value = 1.61
$ws.Cells.Item(5, 18).Value = 1.63
$ws.Cells.Item(5, 19).Value = 2.5
$ws.Cells.Item(5, 21).Value = 2.4
$ws.Cells.Item(5, 24).Value = 24
$ws.Cells.Item(5, 25).Value = 27
$ws.Cells.Item(5, 35).Value = 65
$ws.Cells.Item(5, 41).Value = 65
$ws.Cells.Item(6, 7).Value = 3.4
$ws.Cells.Item(6, 10).Value = 3.8
$ws.Cells.Item(7, 9).Value = 1.72
$ws.Cells.Item(7, 17).Value = 1.84
$ws.Cells.Item(7, 20).Value = 1.85
$ws.Cells.Item(7, 21).Value = 2.1
$ws.Cells.Item(7, 36).Value = 150
$ws.Cells.Item(8, 6).Value = 2.2
$ws.Cells.Item(8, 24).Value = 11
$ws.Cells.Item(8, 26).Value = 28
$ws.Cells.Item(8, 27).Value = 110
$ws.Cells.Item(8, 29).Value = 7.4
$ws.Cells.Item(8, 36).Value = 29
$ws.Cells.Item(9, 8).Value = 2.76
$ws.Cells.Item(10, 6).Value = 1.59
$ws.Cells.Item(10, 11).Value = 4.6
$ws.Cells.Item(10, 17).Value = 1.92
$ws.Cells.Item(11, 10).Value = 3.25
$ws.Cells.Item(12, 6).Value = 1.89
$ws.Cells.Item(12, 7).Value = 2.12
$ws.Cells.Item(12, 8).Value = 4
$ws.Cells.Item(12, 9).Value = 4.4
$ws.Cells.Item(12, 11).Value = 4.1
$ws.Cells.Item(13, 16).Value = 2.68
$ws.Cells.Item(13, 17).Value = 1.47
$ws.Cells.Item(14, 6).Value = 1.35
$ws.Cells.Item(14, 7).Value = 1.53
$ws.Cells.Item(14, 9).Value = 27
$ws.Cells.Item(15, 10).Value = 4.5
$ws.Cells.Item(15, 16).Value = 2.2
$ws.Cells.Item(15, 17).Value = 1.55
$ws.Cells.Item(16, 7).Value = 3
$ws.Cells.Item(16, 8).Value = 2.9
$ws.Cells.Item(16, 10).Value = 2.98
$ws.Cells.Item(16, 16).Value = 1.78
$ws.Cells.Item(19, 17).Value = 1.46
$ws.Cells.Item(20, 10).Value = 2.84
$ws.Cells.Item(20, 17).Value = 2.06
$ws.Cells.Item(23, 8).Value = 4.2
$ws.Cells.Item(23, 11).Value = 3.95
$ws.Cells.Item(24, 7).Value = 2.52
$ws.Cells.Item(24, 8).Value = 3.2
$ws.Cells.Item(24, 16).Value = 1.72
$ws.Cells.Item(24, 17).Value = 1.87
$ws.Cells.Item(25, 4).Value = 'Arsenal'
$ws.Cells.Item(25, 5).Value = 'Aston Villa'
$ws.Cells.Item(25, 6).Value = 1.52
$ws.Cells.Item(25, 7).Value = 1.53
$ws.Cells.Item(25, 8).Value = 7.8
$ws.Cells.Item(25, 9).Value = 8
$ws.Cells.Item(25, 10).Value = 4.6
$ws.Cells.Item(25, 11).Value = 4.7
$ws.Cells.Item(25, 13).Value = 1.07
$ws.Cells.Item(25, 14).Value = 3.9
$ws.Cells.Item(25, 15).Value = 1.32
$ws.Cells.Item(25, 16).Value = 2
$ws.Cells.Item(25, 17).Value = 1.97
$ws.Cells.Item(25, 18).Value = 1.37
$ws.Cells.Item(25, 19).Value = 3.5
$ws.Cells.Item(25, 20).Value = 2.18
$ws.Cells.Item(25, 21).Value = 1.79
$ws.Cells.Item(25, 24).Value = 15
$ws.Cells.Item(25, 25).Value = 22
$ws.Cells.Item(25, 26).Value = 75
$ws.Cells.Item(25, 27).Value = 500
$ws.Cells.Item(25, 28).Value = 7.8
$ws.Cells.Item(25, 29).Value = 10.5
$ws.Cells.Item(25, 30).Value = 30
$ws.Cells.Item(25, 31).Value = 1000
$ws.Cells.Item(25, 32).Value = 8.2
$ws.Cells.Item(25, 34).Value = 30
$ws.Cells.Item(25, 35).Value = 160
$ws.Cells.Item(25, 36).Value = 12.5
$ws.Cells.Item(25, 37).Value = 16.5
$ws.Cells.Item(25, 38).Value = 44
$ws.Cells.Item(25, 39).Value = 1000
$ws.Cells.Item(25, 40).Value = 8.2
$ws.Cells.Item(26, 4).Value = 'Man Utd'
$ws.Cells.Item(26, 5).Value = 'Wolves'
$ws.Cells.Item(26, 6).Value = 1.38
$ws.Cells.Item(26, 7).Value = 1.4
$ws.Cells.Item(26, 8).Value = 9.4
$ws.Cells.Item(26, 9).Value = 9.8
$ws.Cells.Item(26, 10).Value = 5.6
$ws.Cells.Item(26, 11).Value = 5.8
$ws.Cells.Item(26, 13).Value = 1.03
$ws.Cells.Item(26, 14).Value = 0
$ws.Cells.Item(26, 15).Value = 0
$ws.Cells.Item(26, 16).Value = 2.52
$ws.Cells.Item(26, 17).Value = 1.62
$ws.Cells.Item(26, 18).Value = 0
$ws.Cells.Item(26, 19).Value = 0
$ws.Cells.Item(26, 20).Value = 1.94
$ws.Cells.Item(26, 21).Value = 2
$ws.Cells.Item(26, 24).Value = 26
$ws.Cells.Item(26, 25).Value = 36
$ws.Cells.Item(26, 26).Value = 1000
$ws.Cells.Item(26, 27).Value = 390
$ws.Cells.Item(26, 28).Value = 10
$ws.Cells.Item(26, 29).Value = 13
$ws.Cells.Item(26, 30).Value = 36
$ws.Cells.Item(26, 31).Value = 170
$ws.Cells.Item(26, 32).Value = 8.8
$ws.Cells.Item(26, 34).Value = 27
$ws.Cells.Item(26, 35).Value = 1000
$ws.Cells.Item(26, 36).Value = 11.5
$ws.Cells.Item(26, 37).Value = 14
$ws.Cells.Item(26, 38).Value = 34
$ws.Cells.Item(26, 39).Value = 140
$ws.Cells.Item(26, 40).Value = 5.1
$ws.Cells.Item(27, 1).Value = 'Portuguese Segunda Liga'
$ws.Cells.Item(27, 2).Value = '2025-12-30'
$ws.Cells.Item(27, 3).Value = '17:15:00'
$ws.Cells.Item(27, 4).Value = 'Torreense'
$ws.Cells.Item(27, 5).Value = 'Lusitania Futebol Clube'
$ws.Cells.Item(27, 6).Value = 2
$ws.Cells.Item(27, 7).Value = 2.64
$ws.Cells.Item(27, 8).Value = 3.05
$ws.Cells.Item(27, 9).Value = 4.6
$ws.Cells.Item(27, 10).Value = 2.46
$ws.Cells.Item(27, 11).Value = 5.9
$ws.Cells.Item(27, 12).Value = 1.01
$ws.Cells.Item(27, 13).Value = 1.01
$ws.Cells.Item(27, 14).Value = 1.01
$ws.Cells.Item(27, 15).Value = 1.01
$ws.Cells.Item(27, 16).Value = 1.66
$ws.Cells.Item(27, 17).Value = 1.9
$ws.Cells.Item(27, 18).Value = 1.24
$ws.Cells.Item(27, 19).Value = 3.2
$ws.Cells.Item(27, 20).Value = 1.01
$ws.Cells.Item(27, 21).Value = 1.01
$ws.Cells.Item(27, 22).Value = 1.01
$ws.Cells.Item(27, 23).Value = 1.01
$ws.Cells.Item(27, 24).Value = 1000
$ws.Cells.Item(27, 25).Value = 1000
$ws.Cells.Item(27, 26).Value = 1000
$ws.Cells.Item(27, 27).Value = 1000
$ws.Cells.Item(27, 28).Value = 1000
$ws.Cells.Item(27, 29).Value = 1000
$ws.Cells.Item(27, 30).Value = 1000
$ws.Cells.Item(27, 31).Value = 1000
$ws.Cells.Item(27, 32).Value = 1000
$ws.Cells.Item(27, 33).Value = 1000
$ws.Cells.Item(27, 34).Value = 1000
$ws.Cells.Item(27, 35).Value = 1000
$ws.Cells.Item(27, 36).Value = 1000
$ws.Cells.Item(27, 37).Value = 1000
$ws.Cells.Item(27, 38).Value = 1000
$ws.Cells.Item(27, 39).Value = 1000
$ws.Cells.Item(27, 40).Value = 1000
$ws.Cells.Item(27, 41).Value = 1000
